$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.54"
$ws.Range("E2").Value = "'1.01%"
$ws.Range("D3").Value = "'29.41"
$ws.Range("E3").Value = "'-2.60%"
$ws.Range("D4").Value = "'5.150"
$ws.Range("E4").Value = "'0.04%"
$ws.Range("E5").Value = "'1.86%"
$ws.Range("D6").Value = "'6.613"
$ws.Range("E6").Value = "'1.46%"
$ws.Range("D7").Value = "'3.166"
$ws.Range("E7").Value = "'4.93%"
$ws.Range("D8").Value = "'0.8565"
$ws.Range("E8").Value = "'1.91%"
$ws.Range("D9").Value = "'0.8557"
$ws.Range("E9").Value = "'-0.94%"
$ws.Range("D10").Value = "'0.01019"
$ws.Range("E10").Value = "'1,605.91%"
$ws.Range("E11").Value = "'1.92%"
$ws.Range("D12").Value = "'0.07025"
$ws.Range("E12").Value = "'1.44%"
$ws.Range("D13").Value = "'0.03052"
$ws.Range("E13").Value = "'6.77%"
$ws.Range("D14").Value = "'0.09374"
$ws.Range("E14").Value = "'-0.07%"
$ws.Range("D15").Value = "'0.001531"
$ws.Range("E15").Value = "'0.94%"
$ws.Range("D16").Value = "'0.006044"
$ws.Range("E16").Value = "'-0.60%"
$ws.Range("D17").Value = "'3.485"
$ws.Range("E17").Value = "'-0.69%"
$ws.Range("D18").Value = "'2.171"
$ws.Range("E18").Value = "'-2.62%"
$ws.Range("E19").Value = "'1.64%"
$ws.Range("D20").Value = "'0.03328"
$ws.Range("E20").Value = "'1.87%"
$ws.Range("E21").Value = "'-1.07%"
$ws.Range("D22").Value = "'3.320"
$ws.Range("E22").Value = "'-8.02%"
$ws.Range("D23").Value = "'0.04144"
$ws.Range("E23").Value = "'-0.27%"
$ws.Range("D24").Value = "'0.1399"
$ws.Range("E24").Value = "'1.92%"
$ws.Range("E25").Value = "'1.35%"
$ws.Range("D26").Value = "'0.004134"
$ws.Range("E26").Value = "'-4.33%"
$ws.Range("E27").Value = "'2.62%"
$ws.Range("E28").Value = "'3.46%"
$ws.Range("D40").Value = "'0.03726"
$ws.Range("E40").Value = "'0.40%"
$ws.Range("D41").Value = "'0.005889"
$ws.Range("E41").Value = "'10.60%"
$ws.Range("D42").Value = "'0.1069"
$ws.Range("E42").Value = "'1.15%"
$ws.Range("E43").Value = "'-4.28%"
$ws.Range("D44").Value = "'0.008538"
$ws.Range("E44").Value = "'-12.31%"
$ws.Range("D45").Value = "'0.00005285"
$ws.Range("E45").Value = "'3.77%"
$ws.Range("E46").Value = "'0.07%"
$ws.Range("D47").Value = "'0.05799"
$ws.Range("E47").Value = "'-41.96%"
$ws.Range("D48").Value = "'0.002171"
$ws.Range("E48").Value = "'-20.15%"
$ws.Range("E49").Value = "'0.07%"
$ws.Range("E50").Value = "'0.07%"
